$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (b93IZ6K8 / Monterrey vs Club Leon) is replaced with the match that used
# to be row 4 (tA6HgoO8 / Tapatio vs Tepatitlan de Morelos), with refreshed odds.
# Date (column B) is unchanged ("10/11/2024"), so it is left untouched.
$ws.Range("A3").Value = "tA6HgoO8"
$ws.Range("C3").Value = "23:00"
$ws.Range("D3").Value = "MEXICO - LIGA DE EXPANSION MX"
$ws.Range("E3").Value = "Tapatio"
$ws.Range("F3").Value = "Tepatitlan de Morelos"
$ws.Range("G3").Value = 1.87
$ws.Range("H3").Value = 3.25
$ws.Range("I3").Value = 4
$ws.Range("J3").Value = 2.42
$ws.Range("K3").Value = 2.1
$ws.Range("L3").Value = 4.35
$ws.Range("M3").Value = 1.02
$ws.Range("N3").Value = 7.4
$ws.Range("O3").Value = 1.35
$ws.Range("P3").Value = 2.7
$ws.Range("Q3").Value = 2.02
$ws.Range("R3").Value = 1.62
$ws.Range("S3").Value = 1.39
$ws.Range("T3").Value = 2.55
$ws.Range("U3").Value = 1.87
$ws.Range("V3").Value = 1.75
$ws.Range("W3").Value = 6.2
$ws.Range("X3").Value = 8.25
$ws.Range("Y3").Value = 8.5
$ws.Range("Z3").Value = 15.5
$ws.Range("AA3").Value = 16.5
$ws.Range("AB3").Value = 32
$ws.Range("AC3").Value = 8.5
$ws.Range("AD3").Value = 6.4
$ws.Range("AE3").Value = 16.5
$ws.Range("AF3").Value = 90
$ws.Range("AG3").Value = 800
$ws.Range("AH3").Value = 10.25
$ws.Range("AI3").Value = 21
$ws.Range("AJ3").Value = 13.5
$ws.Range("AK3").Value = 65
$ws.Range("AL3").Value = 40
$ws.Range("AM3").Value = 50
$ws.Range("AN3").Value = 3.65
$ws.Range("AO3").Value = 9.25
$ws.Range("AP3").Value = 18.5
$ws.Range("AQ3").Value = 32
$ws.Range("AR3").Value = 70
$ws.Range("AS3").Value = 250
$ws.Range("AT3").Value = 2.5
$ws.Range("AU3").Value = 7.2
$ws.Range("AV3").Value = 65
$ws.Range("AW3").Value = 5.7
$ws.Range("AX3").Value = 22
$ws.Range("AY3").Value = 28
$ws.Range("AZ3").Value = 120
$ws.Range("BA3").Value = 150
$ws.Range("BB3").Value = 350
$ws.Range("BC3").Value = 51
$ws.Range("BD3").Value = 51

# The old row 4 is now redundant (its data lives in row 3) - delete it so the
# sheet shrinks back to A1:BD3.
$ws.Rows("4:4").Delete()
